## "Generate Report for Archive"
##
## The status text "Ready for handoff" moves to "In Translation" everywhere
## it's used (Overview!E2/F2, zh-cn!C2, de-de!C2 all share that string), and
## the now-narrower status columns (Overview E:F, zh-cn/de-de column C) get
## re-sized to fit the new, shorter text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Update every cell that currently reads "Ready for handoff" so the shared
# string collapses onto the new text everywhere it's referenced.
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$zhcn.Range("C2").Value = "In Translation"
$dede.Range("C2").Value = "In Translation"

# Shrink the status columns to match the new (shorter) text.
$overview.Columns("E:F").ColumnWidth = 12.5
$zhcn.Columns("C:C").ColumnWidth = 12.5
$dede.Columns("C:C").ColumnWidth = 12.5
